$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
$c = $cs.Colors(1)
$c.RGB = 255
